# Auto-generated edit script: apply value updates captured in the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 521.26
$ws.Range("I15").Value = 521.26
$ws.Range("K15").Value = 1563.78
$ws.Range("M15").Value = -1394.78
$ws.Range("H61").Value = 93
$ws.Range("I61").Value = 93
$ws.Range("K61").Value = 279
$ws.Range("M61").Value = -107
$ws.Range("H138").Value = 5276.964
$ws.Range("J138").Value = 5419.9614
$ws.Range("L138").Value = 16259.8842
$ws.Range("N138").Value = -26539.8842
$ws.Range("H141").Value = 3550
$ws.Range("I141").Value = 3550
$ws.Range("K141").Value = 10650
$ws.Range("M141").Value = -5470

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 150.33333
$ws.Range("I3").Value = 222.5
$ws.Range("J3").Value = 6
$ws.Range("K3").Value = 222.5
$ws.Range("L3").Value = 6
$ws.Range("M3").Value = -107.5
$ws.Range("N3").Value = -236
$ws.Range("H35").Value = 2365.5
$ws.Range("I35").Value = 2215.5
$ws.Range("J35").Value = 2515.5
$ws.Range("K35").Value = 2215.5
$ws.Range("L35").Value = 2515.5
$ws.Range("M35").Value = -1809.5
$ws.Range("N35").Value = -3327.5
$ws.Range("H61").Value = 1983.2778
$ws.Range("I61").Value = 1592.2307
$ws.Range("K61").Value = 1592.2307
$ws.Range("M61").Value = -1380.2307
$ws.Range("H122").Value = 913635.75
$ws.Range("I122").Value = 2503748.5
$ws.Range("K122").Value = 7511245.5
$ws.Range("M122").Value = -7508795.5
$ws.Range("H132").Value = 1865.6552
$ws.Range("I132").Value = 1781.8889
$ws.Range("K132").Value = 5345.6667
$ws.Range("M132").Value = -2815.6667
$ws.Range("H136").Value = 1983.2778
$ws.Range("I136").Value = 1592.2307
$ws.Range("K136").Value = 4776.6921
$ws.Range("M136").Value = -2226.6921

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 100
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 100
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -380
$ws.Range("H82").Value = 37633
$ws.Range("I82").Value = 6450
$ws.Range("K82").Value = 6450
$ws.Range("M82").Value = -6067
$ws.Range("H85").Value = 37633
$ws.Range("I85").Value = 6450
$ws.Range("K85").Value = 6450
$ws.Range("M85").Value = -5124
$ws.Range("H122").Value = 1979899
$ws.Range("I122").Value = 1979899
$ws.Range("J122").Value = 1979899
$ws.Range("K122").Value = 1979899
$ws.Range("L122").Value = 1979899
$ws.Range("M122").Value = -1974999
$ws.Range("N122").Value = -1989699
$ws.Range("H134").Value = 1436.8422

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6031.5386
$ws.Range("J31").Value = 6532.5454
$ws.Range("L31").Value = 6532.5454
$ws.Range("N31").Value = -7122.5454
$ws.Range("H34").Value = 6031.5386
$ws.Range("J34").Value = 6532.5454
$ws.Range("L34").Value = 6532.5454
$ws.Range("N34").Value = -6936.5454
$ws.Range("H68").Value = 31374.75
$ws.Range("I68").Value = 23500
$ws.Range("J68").Value = 33999.668
$ws.Range("K68").Value = 23500
$ws.Range("L68").Value = 33999.668
$ws.Range("M68").Value = -22751
$ws.Range("N68").Value = -35497.668
$ws.Range("H71").Value = 31374.75
$ws.Range("I71").Value = 23500
$ws.Range("J71").Value = 33999.668
$ws.Range("K71").Value = 70500
$ws.Range("L71").Value = 101999.004
$ws.Range("M71").Value = -66756
$ws.Range("N71").Value = -109487.004
$ws.Range("H88").Value = 1111
$ws.Range("I88").Value = 1111
$ws.Range("K88").Value = 1111
$ws.Range("M88").Value = -705
$ws.Range("H91").Value = 1111
$ws.Range("I91").Value = 1111
$ws.Range("K91").Value = 1111
$ws.Range("M91").Value = 293
$ws.Range("H132").Value = 4035.4375
$ws.Range("I132").Value = 3293.5
$ws.Range("J132").Value = 5272
$ws.Range("K132").Value = 9880.5
$ws.Range("L132").Value = 15816
$ws.Range("M132").Value = -7350.5
$ws.Range("N132").Value = -20876

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 4999.5
$ws.Range("J22").Value = 4999.5
$ws.Range("L22").Value = 14998.5
$ws.Range("N22").Value = -15336.5
$ws.Range("H23").Value = 203.5
$ws.Range("J23").Value = 71.333336
$ws.Range("L23").Value = 214.000008
$ws.Range("N23").Value = -684.000008
$ws.Range("H27").Value = 4999.5
$ws.Range("J27").Value = 4999.5
$ws.Range("L27").Value = 14998.5
$ws.Range("N27").Value = -15202.5
$ws.Range("H35").Value = 345
$ws.Range("I35").Value = 217.5
$ws.Range("J35").Value = 600
$ws.Range("K35").Value = 652.5
$ws.Range("L35").Value = 1800
$ws.Range("M35").Value = -364.5
$ws.Range("N35").Value = -2376
$ws.Range("H38").Value = 139.72
$ws.Range("I38").Value = 60.17647
$ws.Range("K38").Value = 180.52941
$ws.Range("M38").Value = 166.47059
$ws.Range("H50").Value = 254.88889
$ws.Range("I50").Value = 185
$ws.Range("J50").Value = 342.25
$ws.Range("K50").Value = 555
$ws.Range("L50").Value = 1026.75
$ws.Range("M50").Value = -74
$ws.Range("N50").Value = -1988.75
$ws.Range("H53").Value = 254.88889
$ws.Range("I53").Value = 185
$ws.Range("J53").Value = 342.25
$ws.Range("K53").Value = 555
$ws.Range("L53").Value = 1026.75
$ws.Range("M53").Value = -74
$ws.Range("N53").Value = -1988.75
$ws.Range("H63").Value = 1500
$ws.Range("I63").Value = 1500
$ws.Range("K63").Value = 4500
$ws.Range("M63").Value = -3751
$ws.Range("H66").Value = 1500
$ws.Range("I66").Value = 1500
$ws.Range("K66").Value = 13500
$ws.Range("M66").Value = -9756
$ws.Range("H128").Value = 240900
$ws.Range("I128").Value = 240900
$ws.Range("K128").Value = 722700
$ws.Range("M128").Value = -717720

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 395.5
$ws.Range("I55").Value = 348.5
$ws.Range("J55").Value = 560
$ws.Range("K55").Value = 348.5
$ws.Range("L55").Value = 560
$ws.Range("M55").Value = -175.5
$ws.Range("N55").Value = -906
$ws.Range("H122").Value = 3100.7
$ws.Range("I122").Value = 3100.7
$ws.Range("K122").Value = 9302.099999999999
$ws.Range("M122").Value = -6852.099999999999
$ws.Range("H132").Value = 5125.6313
$ws.Range("I132").Value = 5117.1763
$ws.Range("K132").Value = 15351.5289
$ws.Range("M132").Value = -12821.5289

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1617.7812
$ws.Range("J113").Value = 1794
$ws.Range("L113").Value = 5382
$ws.Range("N113").Value = -9722
$ws.Range("H132").Value = 1176.8572
$ws.Range("I132").Value = 1198.037
$ws.Range("K132").Value = 3594.111
$ws.Range("M132").Value = -1064.111

